# Apply the changes described by the commit:
#  - Rename sheet "Salesforce_Case" -> "Sf_Case"
#  - Make "GA_Workbench" the active (selected) tab instead of "Salesforce_Case"
#  - Update the selection / scroll position on both sheets
#  - Add the value "Approved" to cell I2 on the "GA_Workbench" sheet

$wb = $excel.ActiveWorkbook

$wsCase = $wb.Worksheets.Item(1)   # "Salesforce_Case" -> "Sf_Case"
$wsGA   = $wb.Worksheets.Item(2)   # "GA_Workbench"

# Rename the first sheet.
$wsCase.Name = "Sf_Case"

# Add the "Approved" value to the Approved/Rejected column for row 2 of GA_Workbench.
$wsGA.Range("I2").Value = "Approved"

# Update the selection on Sf_Case (was H6, becomes G32) and scroll so F is the
# first visible column (row stays at the top).
$wsCase.Select()
$wsCase.Range("G32").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6

# Make GA_Workbench the active sheet/tab, update its selection (H3) and scroll
# so column D is the first visible column.
$wsGA.Select()
$wsGA.Range("H3").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
